# Apply the data repull/push for dSF (column F) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F5").Value = -3
$ws.Range("F8").Value = -9
$ws.Range("F9").Value = 3
$ws.Range("F14").Value = -4
